$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue 'D2' '29.121.58'
$ws.Range('E2').Value = '  -1.94%  '
Set-TextValue 'D3' '1.851.88'
$ws.Range('E3').Value = '  -0.86%  '
Set-TextValue 'D5' '0.6960'
$ws.Range('E5').Value = '  -4.42%  '
Set-TextValue 'D6' '238.46'
$ws.Range('E6').Value = '  -0.94%  '
$ws.Range('E7').Value = '  +0.07%  '
Set-TextValue 'D8' '0.07649'
$ws.Range('E8').Value = '  +8.04%  '
Set-TextValue 'D9' '0.3030'
$ws.Range('E9').Value = '  -2.90%  '
$ws.Range('E10').Value = '  -4.04%  '
Set-TextValue 'D11' '0.08125'
$ws.Range('E11').Value = '  -1.19%  '
$ws.Range('E12').Value = '  -2.30%  '
Set-TextValue 'D13' '5.222'
$ws.Range('E13').Value = '  -1.59%  '
Set-TextValue 'D14' '1.817.77'
$ws.Range('E14').Value = '  -1.96%  '
Set-TextValue 'D15' '89.09'
$ws.Range('E15').Value = '  -3.25%  '
Set-TextValue 'D16' '29.113.48'
$ws.Range('E16').Value = '  -1.91%  '
$ws.Range('E17').Value = '  -3.48%  '
Set-TextValue 'D18' '13.21'
$ws.Range('E18').Value = '  -0.98%  '
Set-TextValue 'D19' '0.000007727'
$ws.Range('E19').Value = '  -0.62%  '
Set-TextValue 'D20' '236.90'
$ws.Range('E20').Value = '  -4.47%  '
Set-TextValue 'D21' '0.9994'
$ws.Range('E21').Value = '  +0.06%  '
Set-TextValue 'D22' '2.093.84'
$ws.Range('E22').Value = '  -0.23%  '
Set-TextValue 'D23' '1.000'
$ws.Range('E23').Value = '  +0.06%  '
Set-TextValue 'D24' '7.611'
$ws.Range('E24').Value = '  -1.08%  '
Set-TextValue 'D25' '8.982'
$ws.Range('E25').Value = '  -1.80%  '
Set-TextValue 'D26' '161.20'
$ws.Range('E27').Value = '  -5.57%  '
Set-TextValue 'D28' '18.05'
$ws.Range('E28').Value = '  -2.34%  '
Set-TextValue 'D29' '1.989'
$ws.Range('E29').Value = '  -0.73%  '
Set-TextValue 'D30' '1.409'
$ws.Range('E30').Value = '  -1.81%  '
Set-TextValue 'D31' '4.487'
$ws.Range('E31').Value = '  -0.29%  '
Set-TextValue 'D32' '1.488'
$ws.Range('E32').Value = '  -2.13%  '
Set-TextValue 'D33' '4.013'
$ws.Range('E33').Value = '  -3.89%  '
Set-TextValue 'D34' '0.05230'
$ws.Range('E34').Value = '  -0.62%  '
Set-TextValue 'D35' '1.189'
$ws.Range('E35').Value = '  -3.10%  '
$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 'D36' '1.025'
$ws.Range('E36').Value = '  +2.92%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D37' '0.7011'
$ws.Range('E37').Value = '  -6.68%  '
$ws.Range('E38').Value = '  -1.59%  '
Set-TextValue 'D39' '0.01855'
$ws.Range('E39').Value = '  -3.70%  '
Set-TextValue 'D40' '2.678'
$ws.Range('E40').Value = '  -1.99%  '
Set-TextValue 'D41' '0.9331'
$ws.Range('E41').Value = '  +7.92%  '
Set-TextValue 'D42' '6.006'
$ws.Range('E42').Value = '  +0.51%  '
Set-TextValue 'D43' '1.079.31'
$ws.Range('E43').Value = '  +3.33%  '
Set-TextValue 'D44' '0.4265'
$ws.Range('E44').Value = '  -4.22%  '
Set-TextValue 'D45' '70.24'
$ws.Range('E45').Value = '  -0.84%  '
$ws.Range('E46').Value = '  +0.04%  '
Set-TextValue 'D47' '103.01'
$ws.Range('E47').Value = '  -0.87%  '
Set-TextValue 'D48' '1.777'
$ws.Range('E48').Value = '  -1.93%  '
Set-TextValue 'D49' '1.989.12'
$ws.Range('E49').Value = '  -0.60%  '
Set-TextValue 'D50' '9.184'
$ws.Range('E50').Value = '  -3.03%  '
Set-TextValue 'D51' '7.004'
$ws.Range('E51').Value = '  -5.96%  '
